# Update for xlsx.js ver 2.3.0
#
# This reproduces the cell-level effects of the upstream diff:
#   * sheet_1: A1, B1, A2, B2, A3 pick up an explicit (non-default) cell
#     style, while B3 keeps its existing date style but its stored value
#     changes from 41192 (2012-10-10) to 41194.375 (2012-10-12 09:00).
#   * sheet_2 ("シート<2>"): A1, B1, A2, B2 pick up the same explicit style.
#
# The workbook already has an unused "currency" cell style sitting in slot
# 1 of the style table (cellXfs) left over from the template; applying
# that same number format to these ranges makes Excel/the engine reuse
# that existing style slot (style index 1) instead of appending a new
# one, which is what the target file does (same xfId, same cellXfs
# count).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$reusedFormat = "$#,##0_);($#,##0)"

# sheet_1 ----------------------------------------------------------------
$ws1.Range("A1:B2").NumberFormat = $reusedFormat
$ws1.Range("A3").NumberFormat = $reusedFormat

# B3 keeps its existing (date) style - only the underlying serial value
# changes, now carrying a fractional (time-of-day) component.
$ws1.Range("B3").Value2 = 41194.375

# sheet_2 ("シート<2>") ----------------------------------------------------
$ws2.Range("A1:B2").NumberFormat = $reusedFormat
